$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(3, 2).Value  = "0.4.0-snapshot-1"                       # Version
$meta.Cells.Item(6, 2).Value  = "draft"                                   # Status
$meta.Cells.Item(8, 2).Value  = "2024-05-23T12:16:26+00:00"               # Date
$meta.Cells.Item(10, 2).Value = "ANS (https://esante.gouv.fr)"            # Contact

# ---------------------------------------------------------------------------
# 2) Elements sheet: swap the "Mapping: RIM Mapping" (AK) and
#    "Mapping: Spécification métier vers l'extension ROR
#    ContactConfidentialityLevel" (AL) columns - header, widths and data.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

$akCol = 37
$alCol = 38

# -- swap header labels --
$akHeader = $els.Cells.Item(1, $akCol).Value()
$alHeader = $els.Cells.Item(1, $alCol).Value()
$els.Cells.Item(1, $akCol).Value = $alHeader
$els.Cells.Item(1, $alCol).Value = $akHeader

# -- swap the data rows (2 through 6) --
# (skip rows where both cells already hold the same value - e.g. both
#  blank - so we don't turn a "blank" cell into an empty-string cell)
$lastRow = 6
for ($r = 2; $r -le $lastRow; $r++) {
    $akVal = $els.Cells.Item($r, $akCol).Value()
    $alVal = $els.Cells.Item($r, $alCol).Value()
    if ($akVal -ne $alVal) {
        $els.Cells.Item($r, $akCol).Value = $alVal
        $els.Cells.Item($r, $alCol).Value = $akVal
    }
}

# -- swap the column widths to match the (now swapped) content --
# (column AK gets the wide "Spécification métier" width, AL gets the
#  narrower "RIM Mapping" width - i.e. the two stored widths trade places.
#  The ColumnWidth setter here quantises to whole pixels, so the input is
#  chosen to land as close as possible to the target stored widths of
#  83.7734375 / 24.98046875 character-units.)
$els.Columns.Item($akCol).ColumnWidth = 83.0
$els.Columns.Item($alCol).ColumnWidth = 24.1
